$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "E3" = 16.355
    "E4" = 16.431
    "C7" = -12.808
    "A8" = -22.231
    "A10" = -21.754
    "E11" = 16.962
    "A12" = -21.771
    "C14" = -12.773
    "E14" = 16.902
    "C15" = -13.636
    "A18" = -21.572
    "C18" = -10.643
    "E18" = 17.793
    "E19" = 16.452
    "C20" = -11.93
    "E21" = 16.547
    "A25" = -21.832
    "E27" = 16.477
    "C29" = -11.987
    "C30" = -13.347
    "C31" = -13.744
    "E31" = 16.224
    "C35" = -11.986
    "A37" = -20.044
    "E38" = 16.613
    "C40" = -12.782
    "E42" = 16.486
    "C44" = -12.827
    "E44" = 16.732
    "E47" = 16.243
    "C50" = -13.371
    "C54" = -12.661
    "A55" = -21.868
    "E56" = 16.2
    "E58" = 16.602
    "E65" = 17.025
    "A68" = -21.736
    "C68" = -11.001
    "E73" = 16.511
    "C76" = -12.938
    "A77" = -20.586
    "A78" = -20.284
    "A79" = -21.57
    "A80" = -20.098
    "A81" = -21.933
    "A82" = -22.037
    "A84" = -22.169
    "C87" = -13.221
    "C88" = -13.091
    "E90" = 16.328
    "C92" = -11.436
    "E92" = 17.854
    "E94" = 17.828
    "E95" = 17.303
    "C96" = -12.81
    "C98" = -12.888
    "A101" = -21.254
    "C101" = -12.536
    "E101" = 16.717
    "A102" = -19.841
    "C102" = -12.649
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
